# Weekly fruit/vegetable price update: a new weekly record is prepended to
# the data block (row 115), pushing the existing rows 115-138 down to 116-139.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 115, shifting rows 115:138 down to 116:139
# (Excel copies formatting - e.g. the date style on column D - from the
# surrounding rows automatically, same as an interactive "Insert" row).
$ws.Rows(115).Insert()

# Populate the newly inserted row 115 with the new weekly record.
$ws.Range("A115").Value = 7
$ws.Range("B115").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C115").Value = "Ñuble"
$ws.Range("D115").Value = 44508
$ws.Range("E115").Value = 16
$ws.Range("F115").Value = 100112032
$ws.Range("G115").Value = "Zapallo italiano"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 100
$ws.Range("K115").Value = 14000
$ws.Range("L115").Value = 15000
$ws.Range("M115").Value = 14500
$ws.Range("N115").Value = "$/caja 60 unidades"
$ws.Range("O115").Value = "Región del Maule"
$ws.Range("P115").Value = 242
$ws.Range("Q115").Value = 60
$ws.Range("R115").Value = "Hortaliza"
